$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (G=5515)
$ws.Range("H12").Value = 285
$ws.Range("I12").Value = 268.8889
$ws.Range("J12").Value = 430
$ws.Range("K12").Value = 268.8889
$ws.Range("L12").Value = 430
$ws.Range("M12").Value = -98.88889999999998
$ws.Range("N12").Value = -770

# Row 132 (G=44049)
$ws.Range("H132").Value = 6473.4863
$ws.Range("I132").Value = 6844.606
$ws.Range("J132").Value = 3411.75
$ws.Range("K132").Value = 20533.818
$ws.Range("L132").Value = 10235.25
$ws.Range("M132").Value = -18003.818
$ws.Range("N132").Value = -15295.25

# Row 135 (G=44047)
$ws.Range("H135").Value = 912.2093
$ws.Range("I135").Value = 874.5625
$ws.Range("K135").Value = 7871.0625
$ws.Range("M135").Value = -5336.0625

# Row 137 (G=44013)
$ws.Range("H137").Value = 1284778.1
$ws.Range("I137").Value = 2085764.5
$ws.Range("J137").Value = 3200.0667
$ws.Range("K137").Value = 6257293.5
$ws.Range("L137").Value = 9600.2001
$ws.Range("M137").Value = -6254743.5
$ws.Range("N137").Value = -14700.2001

# Row 138 (G=44169)
$ws.Range("H138").Value = 3301.5945
$ws.Range("J138").Value = 3249.7058
$ws.Range("L138").Value = 9749.117400000001
$ws.Range("N138").Value = -20029.1174

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G=27713)
$ws.Range("H2").Value = 821.2
$ws.Range("I2").Value = 755.4167
$ws.Range("J2").Value = 1084.3334
$ws.Range("K2").Value = 755.4167
$ws.Range("L2").Value = 1084.3334
$ws.Range("M2").Value = -642.4167
$ws.Range("N2").Value = -1310.3334

# Row 45 (G=27714)
$ws.Range("H45").Value = 18744.04
$ws.Range("I45").Value = 23357.6
$ws.Range("J45").Value = 3365.5
$ws.Range("K45").Value = 23357.6
$ws.Range("L45").Value = 3365.5
$ws.Range("M45").Value = -22980.6
$ws.Range("N45").Value = -4119.5

# Row 88 (G=12530)
$ws.Range("H88").Value = 3833.1
$ws.Range("I88").Value = 818
$ws.Range("J88").Value = 5843.1665
$ws.Range("K88").Value = 818
$ws.Range("L88").Value = 5843.1665
$ws.Range("M88").Value = -412
$ws.Range("N88").Value = -6655.1665

# Row 91 (G=12530)
$ws.Range("H91").Value = 3833.1
$ws.Range("I91").Value = 818
$ws.Range("J91").Value = 5843.1665
$ws.Range("K91").Value = 818
$ws.Range("L91").Value = 5843.1665
$ws.Range("M91").Value = 586
$ws.Range("N91").Value = -8651.166499999999

# Row 116 (G=27713)
$ws.Range("H116").Value = 821.2
$ws.Range("I116").Value = 755.4167
$ws.Range("J116").Value = 1084.3334
$ws.Range("K116").Value = 755.4167
$ws.Range("L116").Value = 1084.3334
$ws.Range("M116").Value = 1538.5833
$ws.Range("N116").Value = -5672.3334

# Row 123 (G=34107)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G=27713)
$ws.Range("H3").Value = 821.2
$ws.Range("I3").Value = 755.4167
$ws.Range("J3").Value = 1084.3334
$ws.Range("K3").Value = 755.4167
$ws.Range("L3").Value = 1084.3334
$ws.Range("M3").Value = -641.4167
$ws.Range("N3").Value = -1312.3334

# Row 29 (G=2318)
$ws.Range("H29").Value = 5305.3335
$ws.Range("I29").Value = 5305.3335
$ws.Range("K29").Value = 5305.3335
$ws.Range("M29").Value = -5016.3335

# Row 107 (G=27706)
$ws.Range("H107").Value = 1118.5416
$ws.Range("I107").Value = 1033.409
$ws.Range("J107").Value = 2055
$ws.Range("K107").Value = 1033.409
$ws.Range("L107").Value = 2055
$ws.Range("M107").Value = 886.5909999999999
$ws.Range("N107").Value = -5895

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G=5367)
$ws.Range("H22").Value = 1517
$ws.Range("I22").Value = 1523
$ws.Range("J22").Value = 1499
$ws.Range("K22").Value = 1523
$ws.Range("L22").Value = 1499
$ws.Range("M22").Value = -1173
$ws.Range("N22").Value = -2199

# Row 107 (G=27689)
$ws.Range("H107").Value = 1493.5
$ws.Range("I107").Value = 1099
$ws.Range("J107").Value = 1888
$ws.Range("K107").Value = 1099
$ws.Range("L107").Value = 1888
$ws.Range("M107").Value = 821
$ws.Range("N107").Value = -5728

$ws = $wb.Worksheets.Item("CUL")
# Row 121 (G=27878)
$ws.Range("H121").Value = 10310439
$ws.Range("I121").Value = 51000000
$ws.Range("J121").Value = 138049
$ws.Range("K121").Value = 153000000
$ws.Range("L121").Value = 414147
$ws.Range("M121").Value = -152998690
$ws.Range("N121").Value = -416767

$ws = $wb.Worksheets.Item("GSM")
# Row 29 (G=4209)
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

# Row 80 (G=12521)
$ws.Range("H80").Value = 35716484
$ws.Range("I80").Value = 58825696
$ws.Range("J80").Value = 2246.4546
$ws.Range("K80").Value = 58825696
$ws.Range("L80").Value = 2246.4546
$ws.Range("M80").Value = -58824698
$ws.Range("N80").Value = -4242.4546

# Row 83 (G=12521)
$ws.Range("H83").Value = 35716484
$ws.Range("I83").Value = 58825696
$ws.Range("J83").Value = 2246.4546
$ws.Range("K83").Value = 294128480
$ws.Range("L83").Value = 11232.273
$ws.Range("M83").Value = -294123488
$ws.Range("N83").Value = -21216.273

# Row 132 (G=44008)
$ws.Range("H132").Value = 1886.975
$ws.Range("I132").Value = 1403.7273
$ws.Range("J132").Value = 4165.143
$ws.Range("K132").Value = 4211.1819
$ws.Range("L132").Value = 12495.429
$ws.Range("M132").Value = -1681.1819
$ws.Range("N132").Value = -17555.429

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (G=12565)
$ws.Range("H82").Value = 1267.3864
$ws.Range("I82").Value = 1306.2285
$ws.Range("K82").Value = 1306.2285
$ws.Range("M82").Value = -945.2284999999999

# Row 85 (G=12565)
$ws.Range("H85").Value = 1267.3864
$ws.Range("I85").Value = 1306.2285
$ws.Range("K85").Value = 1306.2285
$ws.Range("M85").Value = -58.22849999999994

# Row 136 (G=44060)
$ws.Range("H136").Value = 2446.1316
$ws.Range("I136").Value = 2278.0293
$ws.Range("J136").Value = 3875
$ws.Range("K136").Value = 6834.0879
$ws.Range("L136").Value = 11625
$ws.Range("M136").Value = -4284.0879
$ws.Range("N136").Value = -16725

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (G=27746)
$ws.Range("H107").Value = 1140
$ws.Range("I107").Value = 917
$ws.Range("J107").Value = 1883.3334
$ws.Range("K107").Value = 2751
$ws.Range("L107").Value = 5650.0002
$ws.Range("M107").Value = -831
$ws.Range("N107").Value = -9490.0002

# Row 132 (G=44029)
$ws.Range("H132").Value = 5557646
$ws.Range("I132").Value = 6411992.5
$ws.Range("J132").Value = 4392.875
$ws.Range("K132").Value = 19235977.5
$ws.Range("L132").Value = 13178.625
$ws.Range("M132").Value = -19233447.5
$ws.Range("N132").Value = -18238.625

# Row 136 (G=44031)
$ws.Range("H136").Value = 7104.8887
$ws.Range("I136").Value = 6258.4688
$ws.Range("J136").Value = 13876.25
$ws.Range("K136").Value = 18775.4064
$ws.Range("L136").Value = 41628.75
$ws.Range("M136").Value = -16225.4064
$ws.Range("N136").Value = -46728.75
